$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 4540
$ws1.Range("F4").Value = 482
$ws1.Range("F11").Value = 172
$ws1.Range("F12").Value = 1700
$ws1.Range("F14").Value = 3722
$ws1.Range("F15").Value = 23

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 4540
$ws4.Range("F4").Value = 482
$ws4.Range("F13").Value = 172
$ws4.Range("F16").Value = 1700
$ws4.Range("F18").Value = 3722
$ws4.Range("F19").Value = 23
